$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value = "cost 1"
$ws.Range("L5").Value = "fadil"

$ws.Range("L3").Select()

